{"js": "// P6-2-1a_UART plan: rename the \"P6-2-1a_UART\" token to \"P6-2-1a-UART\"\n// (underscore -> hyphen) everywhere it appears (title + the six\n// \"...captureN.trec\" screen-capture filenames), and relocate the\n// \"_GoBack\" bookmark from its old spot (right before \"in the Shared\n// Files folder\") to its new spot (inside \"Might not be enough time to\n// show this\u2026.?\", right after \"Might not be eno\").\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------\n// 1) Underscore -> hyphen rename, every occurrence.\n//    \"P6-2-1a_UART\"  ->  \"P6-2-1a-UART\"\n// ---------------------------------------------------------------\nconst renameHits = body.search(\"P6-2-1a_UART\", { matchCase: true });\nrenameHits.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < renameHits.items.length; i++) {\n  renameHits.items[i].insertText(\"P6-2-1a-UART\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// ---------------------------------------------------------------\n// 2) Move the \"_GoBack\" bookmark.\n//    a) Remove it from the old location without losing any text: grab\n//       the non-collapsed range \"generated i\" (this literally contains\n//       the bookmark, sitting between \"generated \" and \"in the Shared\n//       Files folder\"), delete it (safe because the range is NOT\n//       collapsed), which also removes the bookmark, then type the\n//       same words back.\n// ---------------------------------------------------------------\nconst oldSpot = body.search(\"generated i\", { matchCase: true });\noldSpot.load(\"items\");\nawait context.sync();\noldSpot.items[0].delete();\nawait context.sync();\n\nconst restoreSpot = body.search(\"been n the Shared\", { matchCase: true });\nrestoreSpot.load(\"items\");\nawait context.sync();\nrestoreSpot.items[0].insertText(\"been generated in the Shared\", Word.InsertLocation.replace);\nawait context.sync();\n\n// b) Re-create \"_GoBack\" at the new location: right after \"Might not be\n//    eno\" (i.e. right before \"ugh time to show this\u2026.?\").\nconst newSpot = body.search(\"Might not be eno\", { matchCase: true });\nnewSpot.load(\"items\");\nawait context.sync();\nconst insertionPoint = newSpot.items[0].getRange(Word.RangeLocation.end);\ninsertionPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# P6-2-1a_UART plan: rename the \"P6-2-1a_UART\" token to \"P6-2-1a-UART\"\n# (underscore -> hyphen) everywhere it appears (title + the six\n# \"...captureN.trec\" screen-capture filenames), and relocate the\n# \"_GoBack\" bookmark from its old spot (right before \"in the Shared\n# Files folder\") to its new spot (inside \"Might not be enough time to\n# show this....?\", right after \"Might not be eno\").\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------\n# 1) Underscore -> hyphen rename, every occurrence.\n#    \"P6-2-1a_UART\"  ->  \"P6-2-1a-UART\"\n# ---------------------------------------------------------------\n$renameRng = $d.Content\n$renameRng.Find.ClearFormatting()\n$renameRng.Find.Replacement.ClearFormatting()\n$renameRng.Find.Text = \"P6-2-1a_UART\"\n$renameRng.Find.Replacement.Text = \"P6-2-1a-UART\"\n$renameRng.Find.Execute(\"P6-2-1a_UART\", $false, $false, $false, $false, $false, $true, 1, $false, \"P6-2-1a-UART\", 2) | Out-Null\n\n# ---------------------------------------------------------------\n# 2) Move the \"_GoBack\" bookmark.\n#    a) Remove it from the old location without losing any text: grab\n#       the non-collapsed range \"generated i\" (this literally contains\n#       the bookmark, sitting between \"generated \" and \"in the Shared\n#       Files folder\"), delete it (safe because the range is NOT\n#       collapsed -- deleting a bookmark's own collapsed Range is a\n#       no-op in Word's object model), which also removes the\n#       bookmark, then type the same words back.\n# ---------------------------------------------------------------\n$oldSpot = $d.Content\n$oldSpot.Find.ClearFormatting()\n$oldSpot.Find.Text = \"generated i\"\n$oldSpot.Find.Execute() | Out-Null\n$oldSpot.Delete()\n\n$restoreSpot = $d.Content\n$restoreSpot.Find.ClearFormatting()\n$restoreSpot.Find.Text = \"been n the Shared\"\n$restoreSpot.Find.Execute() | Out-Null\n$restoreSpot.Text = \"been generated in the Shared\"\n\n# b) Re-create \"_GoBack\" at the new location: right after \"Might not be\n#    eno\" (i.e. right before \"ugh time to show this....?\").\n$newSpot = $d.Content\n$newSpot.Find.ClearFormatting()\n$newSpot.Find.Text = \"Might not be eno\"\n$newSpot.Find.Execute() | Out-Null\n$insertionPoint = $newSpot.Duplicate\n$insertionPoint.Collapse(0)  # wdCollapseEnd\n$d.Bookmarks.Add(\"_GoBack\", $insertionPoint) | Out-Null\n"}
